# "added script for adding high risk country"
# Adds three new columns (HighRiskCountry, TaxOnIncomeFeeType, TaxOnIncomeStatus)
# to the "Institute" worksheet, populates them for every existing institute row
# (plus one brand-new trailing row), and makes "Institute" the active sheet/tab.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Institute")

# Sheets that already contain cells formatted with the header style (bold /
# orange fill / boxed) and the plain boxed "data" style used throughout the
# Institute sheet. We borrow their formatting via copy/paste-special so the
# new cells land on the very same style entries Excel itself would reuse.
$headerStyleSource = $wb.Worksheets.Item("S181372").Range("A1")
$dataStyleSource   = $wb.Worksheets.Item("S205014").Range("A2")

# ---- New header row (row 1) ----
$ws.Range("BD1").Value = "HighRiskCountry"
$ws.Range("BE1").Value = "TaxOnIncomeFeeType"
$ws.Range("BF1").Value = "TaxOnIncomeStatus"

$headerStyleSource.Copy()
$ws.Range("BD1:BF1").PasteSpecial(-4122)

# ---- Existing data rows (2-11) + one new trailing row (12) ----
for ($r = 2; $r -le 12; $r++) {
    $ws.Range("BD$r").Value = "CANADA [124]"
    $ws.Range("BE$r").Value = "Loan Processing Fee [44]"
    $ws.Range("BF$r").Value = "Active [1]"
}

$dataStyleSource.Copy()
$ws.Range("BD2:BF12").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# ---- Make "Institute" the active sheet / tab ----
$ws.Activate()
